# Append 2021-22 NFL odds rows 190-215 (scraped baselines) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowNums = @(190, 191, 192, 193, 194, 195, 196, 197, 198, 199, 200, 201, 202, 203, 204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215)
$Avals = @(1021, 1021, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1024, 1025, 1025)
$Bvals = @(309, 310, 451, 452, 453, 454, 455, 456, 457, 458, 459, 460, 461, 462, 463, 464, 465, 466, 467, 468, 469, 470, 471, 472, 473, 474)
$Cvals = @("V", "H", "V", "H", "V", "H", "V", "H", "V", "H", "V", "H", "V", "H", "V", "H", "V", "H", "V", "H", "V", "H", "V", "H", "V", "H")
$Dvals = @("Denver", "Cleveland", "Cincinnati", "Baltimore", "Carolina", "NYGiants", "Washington", "GreenBay", "KansasCity", "Tennessee", "Atlanta", "Miami", "NYJets", "NewEngland", "Detroit", "LARams", "Philadelphia", "LasVegas", "Chicago", "TampaBay", "Houston", "Arizona", "Indianapolis", "SanFrancisco", "NewOrleans", "Seattle")
$Evals = @(0, 10, 3, 0, 3, 0, 7, 7, 0, 14, 0, 7, 0, 14, 10, 3, 7, 0, 0, 21, 2, 0, 7, 12, 0, 7)
$Fvals = @(0, 0, 10, 10, 0, 5, 0, 7, 0, 13, 13, 0, 7, 17, 6, 14, 0, 17, 3, 14, 3, 17, 6, 0, 10, 0)
$Gvals = @(7, 7, 14, 7, 0, 7, 0, 7, 3, 0, 7, 7, 6, 3, 3, 0, 0, 13, 0, 0, 0, 7, 7, 0, 0, 3)
$Hvals = @(7, 0, 14, 0, 0, 13, 3, 3, 0, 0, 10, 14, 0, 20, 0, 11, 15, 3, 0, 3, 0, 7, 10, 6, 3, 0)
$Ivals = @(14, 17, 41, 17, 3, 25, 10, 24, 3, 27, 30, 28, 13, 54, 19, 28, 22, 33, 3, 38, 5, 31, 30, 18, 13, 10)
$Jvals = @(44.5, 6, 48, 7, 3, 45.5, 50, 7.5, 3, 56, 48, 3, 43.5, 7, 49.5, 13.5, 48.5, 3, 49, 10, 49.5, 14.5, 44.5, 5.5, 3, 44)
$Kvals = @(40, "pk", 46, 7, 2.5, 42.5, 48, 10, 4, 59, "pk", 47.5, 42.5, 7.5, 50.5, 16.5, 48.5, "pk", 47, 13.5, 47.5, 20, 41.5, 3, 5.5, 42)
$Lvals = @(110, -130, 230, -270, -150, 130, 350, -420, -200, 175, -125, 105, 280, -340, 900, -1600, -110, -110, 500, -700, 1200, -3000, 150, -170, -250, 210)
$Mvals = @(2.5, 20, 23.5, 5.5, 0.5, 20.5, 24, 1.5, 7.5, 28, 23.5, 0.5, 19, 2.5, 26, 7.5, 2.5, 24, 20.5, 3.5, 22, 9.5, 20, 1, 3, 20.5)

for ($i = 0; $i -lt $rowNums.Length; $i++) {
    $r = $rowNums[$i]
    $ws.Cells.Item($r, 1).Value = $Avals[$i]
    $ws.Cells.Item($r, 2).Value = $Bvals[$i]
    $ws.Cells.Item($r, 3).Value = $Cvals[$i]
    $ws.Cells.Item($r, 4).Value = $Dvals[$i]
    $ws.Cells.Item($r, 5).Value = $Evals[$i]
    $ws.Cells.Item($r, 6).Value = $Fvals[$i]
    $ws.Cells.Item($r, 7).Value = $Gvals[$i]
    $ws.Cells.Item($r, 8).Value = $Hvals[$i]
    $ws.Cells.Item($r, 9).Value = $Ivals[$i]
    $ws.Cells.Item($r, 10).Value = $Jvals[$i]
    $ws.Cells.Item($r, 11).Value = $Kvals[$i]
    $ws.Cells.Item($r, 12).Value = $Lvals[$i]
    $ws.Cells.Item($r, 13).Value = $Mvals[$i]
}
